# update validasi import data pengguna
# Rename the "program_studi_id" column header to "kode_program_studi"
# in the mahasiswa import template, widen column A slightly, and leave
# the selection where the editor last clicked (D6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C1 held "program_studi_id" - rename it to "kode_program_studi"
$ws.Range("C1").Value = "kode_program_studi"

# Slightly widen column A (nim) to fit the longer values
$ws.Columns("A").ColumnWidth = 11.2

# Leave the cell cursor on D6, matching where editing finished
[void]$ws.Range("D6").Select()
